$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.666.39"
$ws.Range("E2").Value = "  +8.23%  "

$ws.Range("D3").Value = "3.488.12"
$ws.Range("E3").Value = "  +11.56%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "187.96"
$ws.Range("E5").Value = "  +12.85%  "

$ws.Range("D6").Value = "548.04"
$ws.Range("E6").Value = "  +7.47%  "

$ws.Range("D7").Value = "3.481.31"
$ws.Range("E7").Value = "  +11.47%  "

$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +4.65%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "0.630"
$ws.Range("E10").Value = "  +8.06%  "

$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +18.84%  "

$ws.Range("D12").Value = "54.90"
$ws.Range("E12").Value = "  +7.01%  "

$ws.Range("E13").Value = "  +9.59%  "

$ws.Range("D14").Value = "9.34"
$ws.Range("E14").Value = "  +7.95%  "

$ws.Range("D15").Value = "4.049.43"
$ws.Range("E15").Value = "  +11.52%  "

$ws.Range("D16").Value = "3.496.27"
$ws.Range("E16").Value = "  +11.63%  "

$ws.Range("E17").Value = "  +8.02%  "

$ws.Range("D18").Value = "66.650.10"
$ws.Range("E18").Value = "  +8.30%  "

$ws.Range("D19").Value = "18.16"
$ws.Range("E19").Value = "  +8.80%  "

$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  +11.19%  "

$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  +5.95%  "

$ws.Range("D22").Value = "414.95"
$ws.Range("E22").Value = "  +17.06%  "

$ws.Range("D23").Value = "84.90"
$ws.Range("E23").Value = "  +7.91%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "4.24"
$ws.Range("E24").Value = "  +10.72%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "3.90"
$ws.Range("E25").Value = "  +7.80%  "

$ws.Range("D26").Value = "11.09"
$ws.Range("E26").Value = "  +3.85%  "

$ws.Range("D27").Value = "2.91"
$ws.Range("E27").Value = "  +15.29%  "

$ws.Range("D28").Value = "6.12"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "11.81"
$ws.Range("E29").Value = "  +8.75%  "

$ws.Range("D30").Value = "8.74"
$ws.Range("E30").Value = "  +10.86%  "

$ws.Range("D31").Value = "30.06"
$ws.Range("E31").Value = "  +9.53%  "

$ws.Range("D32").Value = "652.02"
$ws.Range("E32").Value = "  +2.62%  "

$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  +7.70%  "

$ws.Range("D34").Value = "11.67"
$ws.Range("E34").Value = "  +6.85%  "

$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +9.84%  "

$ws.Range("D36").Value = "59.18"
$ws.Range("E36").Value = "  +7.18%  "

$ws.Range("D37").Value = "38.45"
$ws.Range("E37").Value = "  +9.44%  "

$ws.Range("D38").Value = "0.0₃0808"
$ws.Range("E38").Value = "  +21.48%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").Value = "0.389"
$ws.Range("E40").Value = "  +7.38%  "

$ws.Range("D41").Value = "0.138"
$ws.Range("E41").Value = "  +15.72%  "

$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  +20.10%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "3.010.85"
$ws.Range("E44").Value = "  +7.66%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.63"
$ws.Range("E45").Value = "  +8.30%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  +17.58%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.27"
$ws.Range("E47").Value = "  +14.11%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0414"
$ws.Range("E48").Value = "  +10.24%  "

$ws.Range("E49").Value = "  +4.36%  "

$ws.Range("D50").Value = "8.85"
$ws.Range("E50").Value = "  +20.53%  "

$ws.Range("E51").Value = "  +8.54%  "
